# fix: connectionTime and disconnectionTime implemented on charging simulation
#
# The charging simulation previously assumed the EV connected at the very
# start of the window. With connectionTime/disconnectionTime now modeled,
# charging only ramps up once the vehicle is actually connected (rows 3-14
# hold flat at the initial SOC / 0 net energy) and the SOC/energy curve is
# recomputed from that later connection point through to the disconnection
# point (rows 15-27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-14: vehicle not yet connected -> SOC stays at the initial 10%,
# no net energy charged, charging power unchanged at 11 kW.
$flatRows = 3..14
foreach ($r in $flatRows) {
    $ws.Cells.Item($r, 2).Value = 11
    $ws.Cells.Item($r, 3).Value = 10
    $ws.Cells.Item($r, 4).Value = 0
}

# Rows 15-26: charging resumes (connected), SOC/energy recalculated with
# the later connection time; charging power stays 11 kW, net energy 2.75 kWh.
$soc = @{
    15 = 16.875
    16 = 23.75
    17 = 30.625
    18 = 37.5
    19 = 44.375
    20 = 51.25
    21 = 58.125
    22 = 65
    23 = 71.875
    24 = 78.75
    25 = 85.625
    26 = 92.5
}
foreach ($r in $soc.Keys) {
    $ws.Cells.Item($r, 2).Value = 11
    $ws.Cells.Item($r, 3).Value = $soc[$r]
    $ws.Cells.Item($r, 4).Value = 2.75
}

# Row 27: battery reaches full charge (100%) and disconnects; charging
# power/energy reflect the tail of the session before disconnection.
$ws.Cells.Item(27, 2).Value = 12
$ws.Cells.Item(27, 3).Value = 100
$ws.Cells.Item(27, 4).Value = 3
